$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("2025-08-05")
$src.Copy($null, $src)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "2025-08-06"

$ws.Range('B2').Value = '転生コロシアム～最弱スキルで最強の女たちを攻略して奴隷ハーレム作ります～'
$ws.Range('C2').Value = 'zunta(作画) はらわたさいぞう(原作)'
$ws.Range('D2').Value = '第31話：完全なる死角①'
$ws.Range('B3').Value = '時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―'
$ws.Range('C3').Value = '光永康則'
$ws.Range('D3').Value = '第６８話『施錠停止』⓵'
$ws.Range('B4').Value = '地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。'
$ws.Range('C4').Value = 'マツモトケンゴ'
$ws.Range('D4').Value = '第６２話　尋問の戦いが始まった（１）'
$ws.Range('B5').Value = '勇者パーティーをクビになったので故郷に帰ったら、メンバー全員がついてきたんだが'
$ws.Range('C5').Value = '絶叫あいす。(漫画) 木の芽(原作) 希(キャラクター原案)'
$ws.Range('D5').Value = '第4話 前編'
$ws.Range('B6').Value = 'このヒーラー、めんどくさい'
$ws.Range('C6').Value = '丹念に発酵(著者)'
$ws.Range('D6').Value = '第89話：盗賊再び'
$ws.Range('B7').Value = '望まぬ不死の冒険者'
$ws.Range('C7').Value = '中曽根ハイジ（漫画） 丘野 優（原作） じゃいあん（キャラクター原案）'
$ws.Range('D7').Value = '第59話　ヴィステルヤ（前編）'
$ws.Range('B8').Value = '剥かせて！竜ケ崎さん'
$ws.Range('C8').Value = '一智和智'
$ws.Range('D8').Value = '大学生編 第13話'
$ws.Range('B9').Value = 'バキ外伝 烈海王は異世界転生しても一向にかまわんッッ'
$ws.Range('C9').Value = '板垣恵介 猪原賽 陸井栄史'
$ws.Range('D9').Value = '第77話　応援(エール)'
$ws.Range('B10').Value = '世界最強の魔女、始めました 〜私だけ『攻略サイト』を見れる世界で自由に生きます〜'
$ws.Range('C10').Value = '戸賀 環 坂木持丸 riritto'
$ws.Range('D10').Value = '第50話②　祝われた家を探索してみた'
$ws.Range('B11').Value = 'バキ外伝　ガイアとシコルスキー　～ときどきノムラ 二人だけど三人暮らし～'
$ws.Range('C11').Value = '板垣恵介 林たかあき'
$ws.Range('D11').Value = '第51話 春の大敵'
$ws.Range('B12').Value = '婚約者に裏切られた錬金術師は、独立して『ざまぁ』します　コミック版'
$ws.Range('C12').Value = '漫画/すたひろ 原作/Y.A'
$ws.Range('D12').Value = 'chapter67【35話②】'
$ws.Range('B13').Value = '魔のものたちは企てる'
$ws.Range('C13').Value = '加藤拓弐(原作) ガしガし(作画)'
$ws.Range('D13').Value = 'コミックス告知'
$ws.Range('B14').Value = 'ハズレ枠の【状態異常スキル】で最強になった俺がすべてを蹂躙するまで'
$ws.Range('C14').Value = '鵜吉しょう（作画） 内々けやき（構成） 篠崎 芳（原作） KWKM（キャラクター原案）'
$ws.Range('D14').Value = '第56話　十河綾香'
$ws.Range('B15').Value = '絶対死なないステラ姫'
$ws.Range('C15').Value = '光永康則 大高稲'
$ws.Range('D15').Value = '第１４話　絶対旅立たない（４）'
$ws.Range('B16').Value = 'ひとりぼっちの異世界攻略'
$ws.Range('C16').Value = 'びび（漫画） 五示正司（原作）'
$ws.Range('D16').Value = '第230話　役立たずの王女'
$ws.Range('B17').Value = '異世界でスローライフを（願望）'
$ws.Range('C17').Value = '長頼（漫画） シゲ（原作） オウカ（キャラクター原案）'
$ws.Range('D17').Value = '第54話　王都一武術大会'
$ws.Range('B18').Value = '彼女を奪ったイケメン美少女がなぜか俺まで狙ってくる'
$ws.Range('C18').Value = '鹿もみじ(漫画) 福田週人(原作) さなだケイスイ(キャラクター原案)'
$ws.Range('D18').Value = '第10話'
$ws.Range('B19').Value = 'クラス転移に巻き込まれたコンビニ店員のおっさん、勇者には必要なかった余り物スキルを駆使して最強となるようです。　コミック版'
$ws.Range('C19').Value = '漫画：結城焔 原作：Narrative Works　日浦あやせ キャラクター原案：鱈'
$ws.Range('D19').Value = 'chapter50【23話③】'
$ws.Range('B20').Value = 'ネットの『推し』とリアルの『推し』が隣に引っ越してきた'
$ws.Range('C20').Value = 'カタケイ（漫画） 遥 透子（原作） 秋乃える（キャラクター原案）'
$ws.Range('D20').Value = '第18話　相合傘'
$ws.Range('B21').Value = '世界最強の騎士は、必ず死ぬヒロインを救うため異世界でも最強の騎士となる〜両手に花を、両手に剣を〜'
$ws.Range('C21').Value = 'KAZU（原作） イソベカズマ（漫画） moino（メカデザイン協力）'
$ws.Range('D21').Value = '第14話（前編）負けイベント'
$ws.Range('B22').Value = '魔王の俺が奴隷エルフを嫁にしたんだが、どう愛でればいい？'
$ws.Range('C22').Value = '原作／手島史詞 キャラクター原案／COMTA 漫画／板垣ハコ'
$ws.Range('D22').Value = '第72話'
$ws.Range('B23').Value = 'サラリーマンが異世界に行ったら四天王になった話'
$ws.Range('C23').Value = '漫画：村光 原作：ベニガシラ'
$ws.Range('D23').Value = '第89話　限界突破'
$ws.Range('B24').Value = '江戸前エルフ'
$ws.Range('C24').Value = '樋口彰彦'
$ws.Range('D24').Value = '#117'
$ws.Range('B25').Value = '元最強探索者のおじさん。美少女配信者を助けて大バズりしてしまった。'
$ws.Range('C25').Value = 'かなたろー(原作) 草壁レイ(漫画)'
$ws.Range('D25').Value = '第4話　美少女、おじさんに「わからせ」られる。①'
$ws.Range('B26').Value = '竜と歩む成り上がり冒険者道～用済みとしてSランクパーティから追放された回復魔術師、捨てられた先で最強の神竜を復活させてしまう ～　コミック版'
$ws.Range('C26').Value = '漫画/＠カリカリうめ 原作/岸本和葉 キャラクター原案/シソ'
$ws.Range('D26').Value = 'chapter50【25話②】'
$ws.Range('B27').Value = '異世界マンチキン　―HP1のままで最強最速ダンジョン攻略―'
$ws.Range('C27').Value = '原作／志瑞 祐 漫画／青桐 良'
$ws.Range('D27').Value = 'ページ128 後始末'
$ws.Range('B28').Value = '序盤で死ぬ最強のサブキャラに転生したので、ゲーム知識で無双する'
$ws.Range('C28').Value = '作画：マエD 原作：新人'
$ws.Range('D28').Value = '第5話(3)'
$ws.Range('B29').Value = 'おっさんはうぜぇぇぇんだよ! ってギルドから追放したくせに、後から復帰要請を出されても遅い。最高の仲間と出会った俺はこっちで最強を目指す!　コミック版'
$ws.Range('C29').Value = '漫画/羽鳥ぴよこ 原作/おうすけ キャラクター原案/エナミカツミ'
$ws.Range('D29').Value = 'chapter37【19話①】'
$ws.Range('B30').Value = '魔王様の街づくり！～最強のダンジョンは近代都市～'
$ws.Range('C30').Value = '吉川英朗（漫画） 月夜 涙（GAノベル/SBクリエイティブ刊）　（原作）'
$ws.Range('D30').Value = '第73話　黒き煤の記憶'
$ws.Range('B31').Value = '絶対に働きたくないダンジョンマスターが惰眠をむさぼるまで'
$ws.Range('C31').Value = '七六（漫画） 鬼影スパナ（原作） よう太（キャラクター原案）'
$ws.Range('D31').Value = '第66話　三つ巴戦、決着'
$ws.Range('B32').Value = '転生したらスライムだった件　クレイマンREVENGE'
$ws.Range('C32').Value = '原作：伏瀬 漫画：カジカ航 キャラクター原案：みっつばー'
$ws.Range('D32').Value = '第30話　略奪'
$ws.Range('B33').Value = 'え、社内システム全てワンオペしている私を解雇ですか？'
$ws.Range('C33').Value = '漫画：伊於 原作：下城米雪 キャラクター原案：icchi'
$ws.Range('D33').Value = '4巻発売告知漫画'
$ws.Range('B34').Value = '四姉妹は夜をおまちかね'
$ws.Range('C34').Value = '保住圭(原作) Bcoca(作画)'
$ws.Range('D34').Value = '第35夜①：『特別』の形はひとつじゃなくて'
$ws.Range('B35').Value = 'アラサーがVTuberになった話。'
$ws.Range('C35').Value = '犬威赤彦(漫画) とくめい(原作) カラスBTK(キャラクター原案)'
$ws.Range('D35').Value = 'コミックス3巻告知イラスト'
$ws.Range('B36').Value = '勇者パーティを追放された【スキルサポーター】、仲間のスキルを解放して最強に成り上がる'
$ws.Range('C36').Value = '作画：なかお 原作：前田氏'
$ws.Range('D36').Value = '第6話(3)'
$ws.Range('B37').Value = 'ハニートラップ・シェアハウス'
$ws.Range('C37').Value = '久慈マサムネ(原作) 神月洸壱(作画)'
$ws.Range('D37').Value = 'コミックス告知'
$ws.Range('B38').Value = '恋人のフリのフリ'
$ws.Range('C38').Value = 'むねひろ(著者)'
$ws.Range('D38').Value = '第8話①'
$ws.Range('B39').Value = '悪役一家の奥方、死に戻りして心を入れ替える。'
$ws.Range('C39').Value = '鏡(漫画) 丘野優(原作) TEDDY(キャラクター原案)'
$ws.Range('D39').Value = '第32話②'
$ws.Range('B40').Value = 'ラーメン大好き小泉さん'
$ws.Range('C40').Value = '鳴見なる'
$ws.Range('D40').Value = '19杯目 名古屋'
$ws.Range('B41').Value = 'ゲーム　オブ　ファミリア-家族戦記-'
$ws.Range('C41').Value = 'Ｄ．Ｐ(作画) 山口ミコト(原作)'
$ws.Range('D41').Value = '第74話①'
$ws.Range('B42').Value = '転生したら没落貴族だったので、【呪言】を極めて家族を救います'
$ws.Range('C42').Value = '作画：アマセケイ 原作：メソポ・たみあ'
$ws.Range('D42').Value = '第6話(3)'
$ws.Range('B43').Value = '最強の少年聖騎士、転生者を狩る'
$ws.Range('C43').Value = '作画：御塩 原作：宇奈木ユラ'
$ws.Range('D43').Value = '第6話(3)'
$ws.Range('B44').Value = '初歩魔法しか使わない謎の老魔法使いが旅をする'
$ws.Range('C44').Value = '山代カゲツ(漫画) やまだのぼる(原作) にじまあるく(キャラクター原案)'
$ws.Range('D44').Value = '第5話②'
$ws.Range('B45').Value = '勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～'
$ws.Range('C45').Value = '漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり'
$ws.Range('D45').Value = '第５１話　英雄を倒す器用貧乏（１）'
$ws.Range('B46').Value = '無能の中の無能王子　スキル【無能】を授かりましたが、周りの女性は【傾国】【傾城】【奸婦】【毒婦】【悪婦】【妖婦】とかです'
$ws.Range('C46').Value = '漫画/一夢 原作/福朗 キャラクター原案/菊池政治'
$ws.Range('D46').Value = 'chapter10【6話①】'
$ws.Range('B47').Value = '宮廷魔導師、追放される　～無能だと追い出された最巧の魔導師は、部下を引き連れて冒険者クランを始めるようです～'
$ws.Range('C47').Value = 'きつね丸（漫画） しんこせい（原作） ろこ（キャラクター原案）'
$ws.Range('D47').Value = '第3話　誰かのための戦い（前編）'
$ws.Range('B48').Value = 'ちゃんと吸えない吸血鬼ちゃん'
$ws.Range('C48').Value = '二式恭介(著者)'
$ws.Range('D48').Value = '第102話：吸血鬼ちゃんと卒業旅行①'
$ws.Range('B49').Value = '傭兵団の料理番'
$ws.Range('C49').Value = '梅木泰祐(漫画) 川井昂(原作) 四季童子(キャラクター原案)'
$ws.Range('D49').Value = '第9話-2'
$ws.Range('B50').Value = '濁る瞳で何を願う ハイセルク戦記'
$ws.Range('C50').Value = 'トルトネン 創-taro 斎藤八呑'
$ws.Range('D50').Value = '第32話 轍'
$ws.Range('B51').Value = '犬と勇者は飾らない'
$ws.Range('C51').Value = '波多ヒロ（漫画） あまなっとう（原作） ヤスダスズヒト（キャラクター原案）'
$ws.Range('D51').Value = '第28.5話　古代妖魔の目覚め'

$wb.Worksheets.Item("Sheet1").Activate()
